$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New team order + updated goal-distribution figures (rows 2-11, columns A-M)
$data = @(
    @("Sūduva",             3, 2, 4, 2, 0, 1, 1, 1, 1, 0, 4, 2),
    @("Hegelmann Litauen",  2, 0, 3, 4, 2, 2, 3, 2, 4, 2, 2, 1),
    @("Kauno Žalgiris",     1, 1, 2, 3, 2, 3, 0, 0, 2, 1, 2, 2),
    @("Banga",              1, 1, 2, 1, 1, 3, 0, 0, 2, 2, 2, 4),
    @("Dainava",            2, 0, 2, 1, 0, 3, 1, 2, 4, 3, 3, 5),
    @("Riteriai",           1, 0, 4, 2, 5, 4, 3, 2, 4, 4, 1, 1),
    @("Džiugas Telšiai",    0, 0, 3, 0, 3, 2, 2, 1, 1, 0, 2, 1),
    @("Šiauliai",           3, 2, 3, 2, 2, 3, 1, 3, 0, 3, 1, 3),
    @("Žalgiris",           2, 5, 1, 1, 2, 2, 2, 0, 4, 2, 1, 3),
    @("Panevėžys",          1, 2, 1, 3, 3, 2, 3, 2, 3, 2, 2, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
